$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing cycle rows (2,3,4) down to (3,4,5) to make room for a
# new, most-recent cycle at row 2. Using Range.Copy(destination) duplicates
# values/types/styles faithfully (shared-string text cells stay text,
# date cells stay dates, numeric cells stay numeric).
$ws.Range("A2:S4").Copy($ws.Range("A3:S5"))

# Duplicate row 5 (now holding the former row-4 data) down into row 6 to
# seed a new row for property 2045374 with matching formatting/values.
$ws.Range("A5:S5").Copy($ws.Range("A6:S6"))

# Row 2: new, most-recent cycle for the same property - one day after what
# is now row 3's "Year Ending" date (12/31/2020 -> 1/1/2021).
$ws.Range("C2").Value = 44197

# Row 6: new property id, and an intentionally blank "Year Ending" date
# (this is the "blank year_ending case" the commit message refers to).
$ws.Range("A6").Value = 2045374
$ws.Range("C6").Value = ""

# Move the active selection to A7, matching the new end-of-data cell.
$ws.Range("A7").Select()
